$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Initial_Screening")

# Set the value "NA" in column G for rows 6, 7, 8, 10, 11 (cells were previously blank)
$ws.Range("G6").Value = "NA"
$ws.Range("G7").Value = "NA"
$ws.Range("G8").Value = "NA"
$ws.Range("G10").Value = "NA"
$ws.Range("G11").Value = "NA"

# Update the active selection to G11
$ws.Range("G11").Select()
